$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "RFA1 - " -> "RFA" + " " + "1 - "  (3 runs, identical rPr)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("RFA1 - ")
if ($found) {
    $s = $rng.Start
    $rng.Text = "RFA 1 - "

    # "RFA" | " " | "1 - "
    $r2 = $d.Range($s + 3, $s + 4)
    $r2.Font.Bold = 0
    $r2.Font.Bold = 1

    $r3 = $d.Range($s + 4, $s + 8)
    $r3.Font.Bold = 0
    $r3.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# Change 2: " visualizza la pagina principale " ->
#           " visualizza la " + "home page"      (2 runs, identical rPr)
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(" visualizza la pagina principale ")
if ($found2) {
    $s2 = $rng2.Start
    $rng2.Text = " visualizza la home page"

    # " visualizza la " | "home page"
    $r2b = $d.Range($s2 + 15, $s2 + 24)
    $r2b.Font.Bold = 0
    $r2b.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# Change 3: " cerca le aziende inserendo la propria città. " ->
#           " cerca le aziende inserendo la propria città."  (trailing space removed)
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute(" cerca le aziende inserendo la propria città. ", $true, $false, $false, $false, $false, $true, 1, $false, " cerca le aziende inserendo la propria città.", 2)
